# Update crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.586.37'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '3.313.36'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.78'
$ws.Range("E5").Value = '  +4.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.91'
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  +3.34%  '
$ws.Range("D9").Value = '3.309.03'
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.177'
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.577'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.18'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  +4.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '637.57'
$ws.Range("E14").Value = '  +11.04%  '
$ws.Range("D15").Value = '3.842.77'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.42'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '67.700.73'
$ws.Range("E17").Value = '  +2.94%  '
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").Value = '3.312.78'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.61'
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.899'
$ws.Range("E22").Value = '  +1.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.56'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.02'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.13'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.00'
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  +3.53%  '
$ws.Range("E28").Value = '  +2.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.69'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.54'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.63'
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '589.62'
$ws.Range("E32").Value = '  +5.68%  '
$ws.Range("D33").Value = '3.925.94'
$ws.Range("E33").Value = '  +5.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.92'
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.53'
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.104'
$ws.Range("E36").Value = '  +1.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.48'
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.128'
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("E40").Value = '  +3.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.68'
$ws.Range("E41").Value = '  +4.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.47'
$ws.Range("E42").Value = '  -1.64%  '
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("D44").Value = '0.0₃0682'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("E47").Value = '  +1.79%  '
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.38'
$ws.Range("E49").Value = '  +12.98%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.53'
$ws.Range("E50").Value = '  +1.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.14'
$ws.Range("E51").Value = '  +4.33%  '
